$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of camera location data for the "Panoche" site (rows 18-25)
$data = @(
    @(18, "Panoche", 1, "High",   "HDPS1C1", 1, 36.69602, -120.79646),
    @(19, "Panoche", 1, "High",   "HDPS1C2", 2, 36.69592, -120.79678),
    @(20, "Panoche", 1, "Medium", "MDPS1C1", 1, 36.69595, -120.7981),
    @(21, "Panoche", 1, "Medium", "MDPS1C2", 2, 36.69595, -120.79779),
    @(22, "Panoche", 1, "Low",    "LDPS1C1", 1, 36.69528, -120.79695),
    @(23, "Panoche", 1, "Low",    "LDPS1C2", 2, 36.69526, -120.79721),
    @(24, "Panoche", 1, "None",   "NDPS1C1", 1, 36.6938,  -120.79295),
    @(25, "Panoche", 1, "None",   "NDPS1C2", 2, 36.69373, -120.79271)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
    $ws.Range("E$r").Value = $row[5]
    $ws.Range("F$r").Value = $row[6]
    $ws.Range("G$r").Value = $row[7]
}

# F24 carries a custom 5-decimal number format
$ws.Range("F24").NumberFormat = "0.00000"

# Update the active selection to match the last edited cell
$null = $ws.Range("F24").Select()
